$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("G2").Value2 = 6263
$ws.Range("J2").Value2 = 5686
$ws.Range("J3").Value2 = 6061
$ws.Range("C4").Value2 = 1835
$ws.Range("E4").Value2 = 2008
$ws.Range("H4").Value2 = 1703
$ws.Range("I4").Value2 = 1773
$ws.Range("J4").Value2 = 1312
$ws.Range("J5").Value2 = 465
$ws.Range("J6").Value2 = 7743
$ws.Range("C7").Value2 = 28379
$ws.Range("E7").Value2 = 26013
$ws.Range("G7").Value2 = 24697
$ws.Range("H7").Value2 = 26014
$ws.Range("I7").Value2 = 26229
$ws.Range("J7").Value2 = 21267

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J7").Value2 = 622
$ws.Range("J8").Value2 = 1339
$ws.Range("J10").Value2 = 146
$ws.Range("J11").Value2 = 329
$ws.Range("J15").Value2 = 234
$ws.Range("J19").Value2 = 621
$ws.Range("J20").Value2 = 442
$ws.Range("J24").Value2 = 66
$ws.Range("J29").Value2 = 1189
$ws.Range("J30").Value2 = 83
$ws.Range("J31").Value2 = 193
$ws.Range("J33").Value2 = 984
$ws.Range("J36").Value2 = 293
$ws.Range("J48").Value2 = 250
$ws.Range("J49").Value2 = 145
$ws.Range("J51").Value2 = 260
$ws.Range("J52").Value2 = 533
$ws.Range("J54").Value2 = 411
$ws.Range("J55").Value2 = 286
$ws.Range("J60").Value2 = 127
$ws.Range("C63").Value2 = 265
$ws.Range("E63").Value2 = 350
$ws.Range("G63").Value2 = 272
$ws.Range("H63").Value2 = 260
$ws.Range("I63").Value2 = 242
$ws.Range("J63").Value2 = 76
$ws.Range("J64").Value2 = 141
$ws.Range("J65").Value2 = 540
$ws.Range("I67").Value2 = 979
$ws.Range("J71").Value2 = 72
$ws.Range("J73").Value2 = 202
$ws.Range("J76").Value2 = 314
$ws.Range("J78").Value2 = 265
$ws.Range("J79").Value2 = 608
$ws.Range("J83").Value2 = 431
$ws.Range("J85").Value2 = 888
$ws.Range("J86").Value2 = 132
$ws.Range("J92").Value2 = 65
$ws.Range("J94").Value2 = 214
$ws.Range("J95").Value2 = 315
$ws.Range("J99").Value2 = 335
$ws.Range("C101").Value2 = 28379
$ws.Range("E101").Value2 = 26013
$ws.Range("G101").Value2 = 24697
$ws.Range("H101").Value2 = 26014
$ws.Range("I101").Value2 = 26229
$ws.Range("J101").Value2 = 21267

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J5").Value2 = 19
$ws.Range("J7").Value2 = 622

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value2 = 100
$ws.Range("J3").Value2 = 65
$ws.Range("J4").Value2 = 22
$ws.Range("J7").Value2 = 329

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value2 = 231
$ws.Range("J3").Value2 = 320
$ws.Range("J6").Value2 = 261
$ws.Range("J7").Value2 = 888

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value2 = 164
$ws.Range("J7").Value2 = 533

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value2 = 369
$ws.Range("J3").Value2 = 407
$ws.Range("J6").Value2 = 451
$ws.Range("J7").Value2 = 1339

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value2 = 128
$ws.Range("J3").Value2 = 159
$ws.Range("J4").Value2 = 13
$ws.Range("J7").Value2 = 431

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value2 = 238
$ws.Range("J3").Value2 = 326
$ws.Range("J6").Value2 = 338
$ws.Range("J7").Value2 = 984

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value2 = 111
$ws.Range("J7").Value2 = 315

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value2 = 155
$ws.Range("J6").Value2 = 190
$ws.Range("J7").Value2 = 540

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value2 = 133
$ws.Range("J6").Value2 = 85
$ws.Range("J7").Value2 = 335

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J6").Value2 = 22
$ws.Range("J7").Value2 = 83

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value2 = 51
$ws.Range("J7").Value2 = 193

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I4").Value2 = 55
$ws.Range("I7").Value2 = 979

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J4").Value2 = 10
$ws.Range("J7").Value2 = 145

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value2 = 197
$ws.Range("J7").Value2 = 411

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value2 = 357
$ws.Range("J3").Value2 = 413
$ws.Range("J6").Value2 = 309
$ws.Range("J7").Value2 = 1189

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value2 = 46
$ws.Range("J6").Value2 = 125
$ws.Range("J7").Value2 = 250

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value2 = 235
$ws.Range("J7").Value2 = 621

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value2 = 47
$ws.Range("J3").Value2 = 66
$ws.Range("J7").Value2 = 314

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J3").Value2 = 29
$ws.Range("J7").Value2 = 146

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value2 = 86
$ws.Range("J7").Value2 = 265

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J2").Value2 = 66
$ws.Range("J3").Value2 = 64
$ws.Range("J7").Value2 = 286

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value2 = 21
$ws.Range("J7").Value2 = 66

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value2 = 171
$ws.Range("J5").Value2 = 18
$ws.Range("J6").Value2 = 176
$ws.Range("J7").Value2 = 608

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value2 = 52
$ws.Range("J7").Value2 = 141

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value2 = 122
$ws.Range("J3").Value2 = 154
$ws.Range("J6").Value2 = 117
$ws.Range("J7").Value2 = 442

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value2 = 95
$ws.Range("J3").Value2 = 95
$ws.Range("J6").Value2 = 88
$ws.Range("J7").Value2 = 293

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value2 = 38
$ws.Range("J7").Value2 = 214

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value2 = 68
$ws.Range("J7").Value2 = 234

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J4").Value2 = 13
$ws.Range("J6").Value2 = 67
$ws.Range("J7").Value2 = 202

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J6").Value2 = 21
$ws.Range("J7").Value2 = 65

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J2").Value2 = 36
$ws.Range("J6").Value2 = 43

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J2").Value2 = 18
$ws.Range("J7").Value2 = 132

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value2 = 62
$ws.Range("J7").Value2 = 260

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value2 = 37
$ws.Range("J7").Value2 = 127

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value2 = 20
$ws.Range("J7").Value2 = 72
